$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (miRNA name, B, C, D, E, F, G)
$newRows = @(
    @("mmu-miR-143-5p",    0, 0, 0, 1, 0, 0),
    @("mmu-miR-133a-5p",   0, 0, 0, 1, 0, 0),
    @("mmu-miR-129-2-3p",  0, 0, 0, 1, 0, 0),
    @("mmu-miR-129-5p",    0, 0, 0, 1, 0, 0),
    @("mmu-miR-434-3p",    0, 0, 0, 1, 1, 1),
    @("mmu-miR-10b-3p",    0, 0, 0, 1, 0, 0),
    @("mmu-miR-27a-5p",    0, 0, 0, 1, 0, 0),
    @("mmu-miR-7036-3p",   0, 0, 0, 1, 0, 0),
    @("mmu-miR-7a-2-3p",   0, 0, 0, 1, 0, 0),
    @("mmu-miR-1a-3p",     0, 0, 0, 1, 0, 0),
    @("mmu-miR-133b-3p",   0, 0, 0, 1, 0, 0),
    @("mmu-miR-668-3p",    0, 0, 0, 0, 0, 1),
    @("mmu-miR-541-5p",    0, 0, 0, 0, 0, 1),
    @("mmu-miR-184-3p",    0, 0, 0, 0, 0, 1),
    @("mmu-miR-6414",      0, 0, 0, 0, 0, 1)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
